$p = $ppt.ActivePresentation

# Slide 2 (1-based) -> ppt/slides/slide2.xml
$s = $p.Slides.Item(2)

# "Group 108" contains "TextBox 109", which holds the
# "{{val:dat_index_count_for_pptx.csv[1:2]}} UPI Defined" field-merge
# placeholder text split across 4 runs.
$grp = $s.Shapes.Item("Group 108")
$tb = $grp.GroupItems.Item("TextBox 109")

# Preserve the box's current (auto-fit) height so that re-writing the run's
# text does not leave the shape's computed layout drifted from its stored
# extent once PowerPoint re-lays the autosize text box out.
$origHeight = $tb.Height

$tr = $tb.TextFrame.TextRange

# The 3rd run in the paragraph holds "[1:2]}}" -> bump the index to "[1:3]}}".
$tr.Runs(3).Text = "[1:3]}}"

# Restore the box height (small epsilon compensates for the Single-precision
# round trip through the Height property so the stored EMU value matches
# the original exactly).
$tb.Height = $origHeight + 0.00005
